$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.04111635684967
$ws.Range("B1").Value = 1.817053198814392
$ws.Range("C1").Value = 2.626778841018677
$ws.Range("D1").Value = 2.414323568344116
$ws.Range("E1").Value = 0.4521180391311646
